$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Fecha" (D) and "Volumen" (M) values per row, reflecting the
# weekly re-shuffle of the underlying rows (other columns for rows 2,3,6,7,8,9,10
# are identical so only D/M change value); rows 4 and 5 additionally swap
# their Q/S/T (unit / price-per-kg / kg-per-unit) details.

$ws.Range("D2").Value = 44322
$ws.Range("M2").Value = 60

$ws.Range("D3").Value = 44313
$ws.Range("M3").Value = 120

$ws.Range("D4").Value = 44327
$ws.Range("M4").Value = 60
$ws.Range("Q4").Value = "$/caja 10 kilos empedrada"
$ws.Range("S4").Value = 11500
$ws.Range("T4").Value = 1

$ws.Range("D5").Value = 44309
$ws.Range("Q5").Value = "$/caja 14 kilos granel"
$ws.Range("S5").Value = 821
$ws.Range("T5").Value = 14

$ws.Range("D6").Value = 44316
$ws.Range("M6").Value = 120

$ws.Range("D7").Value = 44323
$ws.Range("M7").Value = 80

$ws.Range("D8").Value = 44306
$ws.Range("M8").Value = 80

$ws.Range("D9").Value = 44330
$ws.Range("M9").Value = 60

$ws.Range("D10").Value = 44302
